# Append 7 new data rows (8-14) to Sheet1, mirroring the structure of the
# existing rows (2-7): AirportID, Airport, RunwayID, Runway, LDID, LD,
# FunctionID, Function, SegmentID, Segment, CircuitID, Circuit, CCR_ID, CCR_Name.
# LDID (E) and LD (F) are left blank for these new rows, same as in the diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("1","Queenstown","1","05-23","","","1","RWE","1","RWE","10","RWE 1","1","0"),
    @("1","Queenstown","1","05-23","","","1","RWE","1","RWE","20","RWE 2","2","0"),
    @("1","Queenstown","1","05-23","","","1","RWE","1","RWE","30","RWE 3","3","0"),
    @("1","Queenstown","1","05-23","","","2","RCL","2","RCL","40","RCL 1","4","0"),
    @("1","Queenstown","1","05-23","","","2","RCL","2","RCL","50","RCL 2","5","0"),
    @("1","Queenstown","1","05-23","","","2","RCL","2","RCL","60","RCL 3","6","0"),
    @("1","Queenstown","1","05-23","","","7","TWY","7","TWY","110","TWY","11","0")
)

$rowIndex = 8
foreach ($row in $data) {
    $colIndex = 1
    foreach ($val in $row) {
        if ($val -ne "") {
            if ($val -match '^-?\d+$') {
                # Numeric-looking values are stored as text in this sheet
                # (matching the existing rows), so force text with a
                # leading apostrophe just like typing it in Excel would.
                $ws.Cells.Item($rowIndex, $colIndex).Value = "'" + $val
            } else {
                $ws.Cells.Item($rowIndex, $colIndex).Value = $val
            }
        }
        $colIndex++
    }
    $rowIndex++
}
